# Updated cryptos list on Sun Mar 24 20:08:05 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the coinranking
# snapshot, plus a handful of rank-reshuffle swaps where the coin at a given
# row changed (B/C/D/E all updated together for those rows).
#
# Note: several "Price" strings look like plain decimals (e.g. "563.22"),
# which Excel's Range.Value setter would silently reinterpret as a number
# (losing the intended text type / trailing zeros, e.g. "1.00" -> 1). Those
# assignments are prefixed with a leading apostrophe to force plain-text
# entry, exactly like typing '563.22 into a cell in the Excel UI; the
# apostrophe itself is not stored as part of the cell text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.810.65"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "3.379.49"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'563.22"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").Value = "'175.77"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("D8").Value = "3.375.40"
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("E10").Value = "  +2.55%  "

$ws.Range("D11").Value = "'0.633"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "'53.75"

$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = "  -1.12%  "

$ws.Range("D14").Value = "'9.25"
$ws.Range("E14").Value = "  +1.19%  "

$ws.Range("D15").Value = "3.929.72"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "'18.17"
$ws.Range("E17").Value = "  -1.01%  "

$ws.Range("D18").Value = "3.378.57"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "65.878.75"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").Value = "'11.87"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").Value = "'463.03"
$ws.Range("E22").Value = "  -2.56%  "

$ws.Range("D23").Value = "'4.92"
$ws.Range("E23").Value = "  -0.81%  "

$ws.Range("D24").Value = "'14.98"
$ws.Range("E24").Value = "  +10.88%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'89.69"
$ws.Range("E25").Value = "  +2.94%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'4.10"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").Value = "'10.65"

$ws.Range("D29").Value = "'8.71"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").Value = "'31.04"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").Value = "'6.61"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("D32").Value = "'11.45"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("D33").Value = "'581.26"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("D34").Value = "'62.24"
$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("D38").Value = "'3.56"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").Value = "'36.01"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").Value = "'0.378"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").Value = "0.0₃0746"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("D42").Value = "3.101.04"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("D43").Value = "'2.83"
$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.134"
$ws.Range("E45").Value = "  -1.04%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.45"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.17"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").Value = "'141.29"
$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D50").Value = "'2.56"
$ws.Range("E50").Value = "  +8.94%  "

$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'2.56"
$ws.Range("E51").Value = "  -1.98%  "

